# TestSet30.xlsx - "Scripts" sheet: add randomised id values to column A
# for rows 10-16 (the F:I "Infer" formulas already append $A<row>, so they
# recalc automatically), widen column E, and move the selection to A17.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scripts")

# New data in column A (previously blank) for rows 10-16.
$ws.Range("A10").Value = 7652
$ws.Range("A11").Value = 8384
$ws.Range("A12").Value = 6600
$ws.Range("A13").Value = 2220
$ws.Range("A14").Value = 6612
$ws.Range("A15").Value = 9492
$ws.Range("A16").Value = 5756

# Widen column E (best-fit-style widen to fit the long formula results).
$ws.Columns.Item(5).ColumnWidth = 103.71

# Update the view: scroll back to the top and move the selection to A17.
$ws.Range("A17").Select()

$wb.Save()
